$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated TPM pipeline no longer produces any "Resolving-Mac" target-cluster
# rows, so the shared string "Resolving-Mac" and its rows (originally rows 5, 9
# and 13, one per sending cluster) are removed. Deleting from the bottom up keeps
# the remaining row numbers stable while each delete is applied.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(5).Delete()

# The remaining 9 rows (ECs/FAPs/MuSCs x ECs/FAPs/MuSCs) keep their sending and
# target cluster labels, but all of the numeric statistics (columns E:T) were
# recomputed with the new TPM values. Write the refreshed figures in place.
$newValues = @{
    "G2" = 0.3227736666666667
    "H2" = 0.968321
    "I2" = 0.1416094457286952
    "J2" = 0.1416094457286952
    "M2" = 0.2282403333333334
    "N2" = 0.684721
    "O2" = 0.1514399067192403
    "P2" = 0.1514399067192403
    "Q2" = 0.07366996927122223
    "R2" = 0.663029723441
    "S2" = 0.02144532125171692
    "T2" = 0.02144532125171691
    "G3" = 0.3227736666666667
    "H3" = 0.968321
    "I3" = 0.1416094457286952
    "J3" = 0.1416094457286952
    "O3" = 0.763937263734893
    "P3" = 0.763937263734893
    "Q3" = 0.3716275053498889
    "R3" = 3.344647548149001
    "S3" = 0.1081807324889942
    "T3" = 0.1081807324889942
    "G4" = 0.3227736666666667
    "H4" = 0.968321
    "I4" = 0.1416094457286952
    "J4" = 0.1416094457286952
    "M4" = 0.127538
    "N4" = 0.3826140000000001
    "O4" = 0.08462282954586674
    "P4" = 0.08462282954586674
    "Q4" = 0.04116590789933334
    "R4" = 0.370493171094
    "S4" = 0.01198339198798404
    "T4" = 0.01198339198798404
    "I5" = 0.8226066833587575
    "J5" = 0.8226066833587576
    "M5" = 0.2282403333333334
    "N5" = 0.684721
    "O5" = 0.1514399067192403
    "P5" = 0.1514399067192403
    "Q5" = 0.4279475057154445
    "R5" = 3.851527551439001
    "S5" = 0.1245754793944739
    "T5" = 0.1245754793944739
    "I6" = 0.8226066833587575
    "J6" = 0.8226066833587576
    "O6" = 0.763937263734893
    "P6" = 0.763937263734893
    "S6" = 0.6284198988151247
    "T6" = 0.6284198988151248
    "I7" = 0.8226066833587575
    "J7" = 0.8226066833587576
    "M7" = 0.127538
    "N7" = 0.3826140000000001
    "O7" = 0.08462282954586674
    "P7" = 0.08462282954586674
    "Q7" = 0.2391320069806667
    "R7" = 2.152188062826001
    "S7" = 0.0696113051491589
    "T7" = 0.06961130514915892
    "E8" = 2
    "F8" = 0.6666666666666666
    "G8" = 0.081563
    "H8" = 0.244689
    "I8" = 0.03578387091254728
    "J8" = 0.03578387091254728
    "M8" = 0.2282403333333334
    "N8" = 0.684721
    "O8" = 0.1514399067192403
    "P8" = 0.1514399067192403
    "Q8" = 0.01861596630766667
    "R8" = 0.167543696769
    "S8" = 0.005419106073049496
    "T8" = 0.005419106073049495
    "E9" = 2
    "F9" = 0.6666666666666666
    "G9" = 0.081563
    "H9" = 0.244689
    "I9" = 0.03578387091254728
    "J9" = 0.03578387091254728
    "O9" = 0.763937263734893
    "P9" = 0.763937263734893
    "Q9" = 0.09390807661566666
    "R9" = 0.8451726895410001
    "S9" = 0.027336632430774
    "T9" = 0.027336632430774
    "E10" = 2
    "F10" = 0.6666666666666666
    "G10" = 0.081563
    "H10" = 0.244689
    "I10" = 0.03578387091254728
    "J10" = 0.03578387091254728
    "M10" = 0.127538
    "N10" = 0.3826140000000001
    "O10" = 0.08462282954586674
    "P10" = 0.08462282954586674
    "Q10" = 0.010402381894
    "R10" = 0.09362143704600001
    "S10" = 0.003028132408723787
    "T10" = 0.003028132408723787
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
